$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 2000
$ws.Range("I7").Value = 2000
$ws.Range("K7").Value = 2000
$ws.Range("M7").Value = -1888

$ws.Range("H14").Value = 2000
$ws.Range("I14").Value = 2000
$ws.Range("K14").Value = 2000
$ws.Range("M14").Value = -1809

$ws.Range("H33").Value = 266.41177
$ws.Range("I33").Value = 220.5625
$ws.Range("K33").Value = 220.5625
$ws.Range("M33").Value = 8.4375

$ws.Range("H98").Value = 2585.5715
$ws.Range("I98").Value = 2799.842
$ws.Range("J98").Value = 550
$ws.Range("K98").Value = 2799.842
$ws.Range("L98").Value = 550
$ws.Range("M98").Value = -1301.842
$ws.Range("N98").Value = -3546

$ws.Range("H122").Value = 2585.5715
$ws.Range("I122").Value = 2799.842
$ws.Range("J122").Value = 550
$ws.Range("K122").Value = 8399.526
$ws.Range("L122").Value = 1650
$ws.Range("M122").Value = -5949.526
$ws.Range("N122").Value = -6550

$ws.Range("H125").Value = 1240.7142
$ws.Range("I125").Value = 950
$ws.Range("J125").Value = 1289.1666
$ws.Range("K125").Value = 8550
$ws.Range("L125").Value = 11602.4994
$ws.Range("M125").Value = -6090
$ws.Range("N125").Value = -16522.4994

$ws.Range("H131").Value = 1948.3334
$ws.Range("I131").Value = 1839
$ws.Range("J131").Value = 2495
$ws.Range("K131").Value = 5517
$ws.Range("L131").Value = 7485
$ws.Range("M131").Value = -477
$ws.Range("N131").Value = -17565

$ws.Range("H141").Value = 1455
$ws.Range("I141").Value = 1132.5
$ws.Range("J141").Value = 2100
$ws.Range("K141").Value = 3397.5
$ws.Range("L141").Value = 6300
$ws.Range("M141").Value = 1782.5
$ws.Range("N141").Value = -16660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 10000000
$ws.Range("I11").Value = 10000000
$ws.Range("K11").Value = 10000000
$ws.Range("M11").Value = -9999856

$ws.Range("H17").Value = 4500
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 4500
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4500
$ws.Range("N17").Value = -4846
$ws.Range("M17").ClearContents()

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H61").Value = 2782.9033
$ws.Range("I61").Value = 1922.9565
$ws.Range("J61").Value = 5255.25
$ws.Range("K61").Value = 1922.9565
$ws.Range("L61").Value = 5255.25
$ws.Range("M61").Value = -1710.9565
$ws.Range("N61").Value = -5679.25

$ws.Range("H110").Value = 1040.1818
$ws.Range("I110").Value = 898.1053000000001
$ws.Range("J110").Value = 1940
$ws.Range("K110").Value = 898.1053000000001
$ws.Range("L110").Value = 1940
$ws.Range("M110").Value = 1146.8947
$ws.Range("N110").Value = -6030

$ws.Range("H122").Value = 2681.5
$ws.Range("I122").Value = 1541.3334
$ws.Range("J122").Value = 3365.6
$ws.Range("K122").Value = 4624.0002
$ws.Range("L122").Value = 10096.8
$ws.Range("M122").Value = -2174.0002
$ws.Range("N122").Value = -14996.8

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 2782.9033
$ws.Range("I136").Value = 1922.9565
$ws.Range("J136").Value = 5255.25
$ws.Range("K136").Value = 5768.8695
$ws.Range("L136").Value = 15765.75
$ws.Range("M136").Value = -3218.8695
$ws.Range("N136").Value = -20865.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 497.4
$ws.Range("I8").Value = 497.4
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 497.4
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -357.4
$ws.Range("N8").ClearContents()

$ws.Range("H10").Value = 1000
$ws.Range("I10").Value = 1000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -860
$ws.Range("N10").ClearContents()

$ws.Range("H64").Value = 1495.7778
$ws.Range("I64").Value = 3277.25
$ws.Range("J64").Value = 745.6842
$ws.Range("K64").Value = 3277.25
$ws.Range("L64").Value = 745.6842
$ws.Range("M64").Value = -3052.25
$ws.Range("N64").Value = -1195.6842

$ws.Range("H67").Value = 1495.7778
$ws.Range("I67").Value = 3277.25
$ws.Range("J67").Value = 745.6842
$ws.Range("K67").Value = 3277.25
$ws.Range("L67").Value = 745.6842
$ws.Range("M67").Value = -2497.25
$ws.Range("N67").Value = -2305.6842

$ws.Range("H98").Value = 59363.668
$ws.Range("J98").Value = 59363.668
$ws.Range("L98").Value = 59363.668
$ws.Range("N98").Value = -65353.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 150
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 300
$ws.Range("L7").Value = 600
$ws.Range("M7").Value = -188
$ws.Range("N7").Value = -824

$ws.Range("H9").Value = 1119.8182
$ws.Range("I9").Value = 620
$ws.Range("J9").Value = 1307.25
$ws.Range("K9").Value = 1860
$ws.Range("L9").Value = 3921.75
$ws.Range("M9").Value = -1636
$ws.Range("N9").Value = -4369.75

$ws.Range("H13").Value = 50
$ws.Range("I13").Value = 55
$ws.Range("K13").Value = 165
$ws.Range("M13").Value = 3

$ws.Range("H92").Value = 679
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 708.8333
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 2126.4999
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -4622.4999

$ws.Range("H105").Value = 908000000
$ws.Range("J105").Value = 908000000
$ws.Range("L105").Value = 2724000000
$ws.Range("N105").Value = -2724005242

$ws.Range("H131").Value = 1192401.5
$ws.Range("I131").Value = 2912.0527
$ws.Range("J131").Value = 1540098.4
$ws.Range("K131").Value = 8736.158100000001
$ws.Range("L131").Value = 4620295.199999999
$ws.Range("M131").Value = -3696.158100000001
$ws.Range("N131").Value = -4630375.199999999

$ws.Range("H137").Value = 40283.62
$ws.Range("I137").Value = 2257.375
$ws.Range("K137").Value = 6772.125
$ws.Range("M137").Value = -1672.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()

$ws.Range("H113").Value = 1177.6428
$ws.Range("I113").Value = 757.125
$ws.Range("J113").Value = 1738.3334
$ws.Range("K113").Value = 757.125
$ws.Range("L113").Value = 1738.3334
$ws.Range("M113").Value = 1412.875
$ws.Range("N113").Value = -6078.3334

$ws.Range("H122").Value = 3687.4285
$ws.Range("I122").Value = 3924.889
$ws.Range("K122").Value = 11774.667
$ws.Range("M122").Value = -9324.667000000001

$ws.Range("H132").Value = 2387.7576
$ws.Range("I132").Value = 2041.579
$ws.Range("J132").Value = 2857.5715
$ws.Range("K132").Value = 6124.737
$ws.Range("L132").Value = 8572.7145
$ws.Range("M132").Value = -3594.737
$ws.Range("N132").Value = -13632.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 100000
$ws.Range("I14").Value = 100000
$ws.Range("K14").Value = 100000
$ws.Range("M14").Value = -99828

$ws.Range("H62").Value = 11000
$ws.Range("I62").Value = 12000
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 12000
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -11376
$ws.Range("N62").Value = -11248

$ws.Range("H65").Value = 11000
$ws.Range("I65").Value = 12000
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 36000
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -32880
$ws.Range("N65").Value = -36240

$ws.Range("H68").Value = 1835.6428
$ws.Range("I68").Value = 1818.1818
$ws.Range("J68").Value = 1899.6666
$ws.Range("K68").Value = 1818.1818
$ws.Range("L68").Value = 1899.6666
$ws.Range("M68").Value = -1069.1818
$ws.Range("N68").Value = -3397.6666

$ws.Range("H71").Value = 1835.6428
$ws.Range("I71").Value = 1818.1818
$ws.Range("J71").Value = 1899.6666
$ws.Range("K71").Value = 9090.909
$ws.Range("L71").Value = 9498.333000000001
$ws.Range("M71").Value = -5346.909
$ws.Range("N71").Value = -16986.333

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()

$ws.Range("H75").Value = 55665.5
$ws.Range("J75").Value = 55665.5
$ws.Range("L75").Value = 55665.5
$ws.Range("N75").Value = -57537.5

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()

$ws.Range("H78").Value = 55665.5
$ws.Range("J78").Value = 55665.5
$ws.Range("L78").Value = 166996.5
$ws.Range("N78").Value = -176356.5

$ws.Range("H80").Value = 15000
$ws.Range("J80").Value = 15000
$ws.Range("L80").Value = 15000
$ws.Range("N80").Value = -17246

$ws.Range("H83").Value = 15000
$ws.Range("J83").Value = 15000
$ws.Range("L83").Value = 45000
$ws.Range("N83").Value = -56232

$ws.Range("H122").Value = 5936.4688
$ws.Range("I122").Value = 6880.5
$ws.Range("J122").Value = 4363.0835
$ws.Range("K122").Value = 20641.5
$ws.Range("L122").Value = 13089.2505
$ws.Range("M122").Value = -18191.5
$ws.Range("N122").Value = -17989.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 34826.75
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 34826.75
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 34826.75
$ws.Range("N12").Value = -35110.75
$ws.Range("M12").ClearContents()

Write-Output "Applied all changes"